$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (old "Peter/Parker" and shift old row4 "Liat/Mulian" data into row2)
# Target final data:
# Row1: A1=208063511 B1=Omer C1=Avisror D1=12345 E1=TRUE
# Row2: A2=313301129 B2=Liat C2=Mulian D2="0546500" (text) E2=TRUE

$ws.Range("A2").Value = 313301129
$ws.Range("B2").Value = "Liat"
$ws.Range("C2").Value = "Mulian"
$ws.Range("E2").Value = $true

# D2 should be text "0546500" (preserve leading zero) - set number format to text first
$ws.Columns("A").NumberFormat = "@"
$ws.Columns("D").NumberFormat = "@"
$ws.Columns("D").ColumnWidth = 9.140625
$ws.Range("D2").Value = "0546500"

# Remove rows 3 and 4 entirely
$ws.Rows("3:4").Delete()

# Update selection to D10 as per sheetView change
$ws.Range("D10").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
